# Apply "added harvard case classification" edit:
# - Insert a new average_doctor column, shifting the previous one to average_doctor_old
# - Recompute the average/variance/stdDev statistics for every "_old" classifier column
#   and for the average_doctor / average_doctor_old summary columns (rows 4-13).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header labels - average_doctor becomes average_doctor_old, and vice versa
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4: updated stats
$ws.Range("E4").Value = 0.667
$ws.Range("F4").Value = 0.062
$ws.Range("G4").Value = 0.249
$ws.Range("N4").Value = 0.611
$ws.Range("O4").Value = 0.08
$ws.Range("P4").Value = 0.283
$ws.Range("Q4").Value = 0.444
$ws.Range("R4").Value = 0.173
$ws.Range("S4").Value = 0.416
$ws.Range("W4").Value = 0.556
$ws.Range("X4").Value = 0.006
$ws.Range("Y4").Value = 0.079
$ws.Range("AI4").Value = 0.611
$ws.Range("AJ4").Value = 0.08
$ws.Range("AK4").Value = 0.283
$ws.Range("AU4").Value = 0.41
$ws.Range("AV4").Value = 0
$ws.Range("AW4").Value = 0.013
$ws.Range("BA4").Value = 1.611
$ws.Range("BB4").Value = 0
$ws.Range("BC4").Value = 0
$ws.Range("BG4").Value = 0.667
$ws.Range("BH4").Value = 0.222
$ws.Range("BI4").Value = 0.471
$ws.Range("BM4").Value = 0.444
$ws.Range("BN4").Value = 0.006
$ws.Range("BO4").Value = 0.079
$ws.Range("BP4").Value = 0.537
$ws.Range("BQ4").Value = 0.714

# Row 5: updated stats
$ws.Range("E5").Value = 0.8
$ws.Range("F5").Value = 0.08
$ws.Range("G5").Value = 0.283
$ws.Range("N5").Value = 0.756
$ws.Range("O5").Value = 0.031
$ws.Range("P5").Value = 0.175
$ws.Range("Q5").Value = 0.178
$ws.Range("R5").Value = 0.019
$ws.Range("S5").Value = 0.137
$ws.Range("W5").Value = 0.511
$ws.Range("X5").Value = 0.122
$ws.Range("Y5").Value = 0.35
$ws.Range("AI5").Value = 0.578
$ws.Range("AJ5").Value = 0.09
$ws.Range("AK5").Value = 0.3
$ws.Range("AU5").Value = 0.756
$ws.Range("AV5").Value = 0.031
$ws.Range("AW5").Value = 0.175
$ws.Range("BA5").Value = 1.177
$ws.Range("BB5").Value = 0.005
$ws.Range("BC5").Value = 0.068
$ws.Range("BG5").Value = 0.422
$ws.Range("BH5").Value = 0.09
$ws.Range("BI5").Value = 0.3
$ws.Range("BM5").Value = 0.344
$ws.Range("BN5").Value = 0.015
$ws.Range("BO5").Value = 0.123
$ws.Range("BP5").Value = 0.392
$ws.Range("BQ5").Value = 0.426

# Row 6: updated stats
$ws.Range("E6").Value = 0.727
$ws.Range("N6").Value = 0.676
$ws.Range("Q6").Value = 0.254
$ws.Range("W6").Value = 0.533
$ws.Range("AI6").Value = 0.594
$ws.Range("AU6").Value = 0.532
$ws.Range("BA6").Value = 1.356
$ws.Range("BG6").Value = 0.517
$ws.Range("BM6").Value = 0.388
$ws.Range("BP6").Value = 0.452
$ws.Range("BQ6").Value = 0.529

# Row 7: updated stats
$ws.Range("E7").Value = 0.769
$ws.Range("N7").Value = 0.722
$ws.Range("Q7").Value = 0.202
$ws.Range("W7").Value = 0.519
$ws.Range("AI7").Value = 0.584
$ws.Range("AU7").Value = 0.647
$ws.Range("BA7").Value = 1.241
$ws.Range("BG7").Value = 0.455
$ws.Range("BM7").Value = 0.36
$ws.Range("BP7").Value = 0.414
$ws.Range("BQ7").Value = 0.462

# Row 8: updated stats
$ws.Range("E8").Value = 0.909
$ws.Range("F8").Value = 0.008999999999999999
$ws.Range("G8").Value = 0.095
$ws.Range("N8").Value = 0.9370000000000001
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0.012
$ws.Range("Q8").Value = 0.262
$ws.Range("R8").Value = 0.089
$ws.Range("S8").Value = 0.298
$ws.Range("W8").Value = 0.74
$ws.Range("X8").Value = 0.002
$ws.Range("Y8").Value = 0.048
$ws.Range("AI8").Value = 0.828
$ws.Range("AJ8").Value = 0.008999999999999999
$ws.Range("AK8").Value = 0.097
$ws.Range("AU8").Value = 0.681
$ws.Range("AV8").Value = 0.026
$ws.Range("AW8").Value = 0.161
$ws.Range("BA8").Value = 1.634
$ws.Range("BB8").Value = 0.037
$ws.Range("BC8").Value = 0.193
$ws.Range("BG8").Value = 0.584
$ws.Range("BH8").Value = 0.173
$ws.Range("BI8").Value = 0.416
$ws.Range("BM8").Value = 0.537
$ws.Range("BN8").Value = 0.096
$ws.Range("BO8").Value = 0.311
$ws.Range("BP8").Value = 0.545
$ws.Range("BQ8").Value = 0.581

# Row 9: updated stats
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("AI9").Value = 1
$ws.Range("AJ9").Value = 0
$ws.Range("AK9").Value = 0
$ws.Range("BA9").Value = 1.667
$ws.Range("BM9").Value = 0.667
$ws.Range("BP9").Value = 0.556
$ws.Range("BQ9").Value = 0.5669999999999999

# Row 10: updated stats
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("W10").Value = 1
$ws.Range("X10").Value = 0
$ws.Range("Y10").Value = 0
$ws.Range("AI10").Value = 1
$ws.Range("AJ10").Value = 0
$ws.Range("AK10").Value = 0
$ws.Range("BA10").Value = 2.001
$ws.Range("BM10").Value = 0.667
$ws.Range("BP10").Value = 0.667
$ws.Range("BQ10").Value = 0.7

# Row 11: updated stats
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("W11").Value = 1
$ws.Range("X11").Value = 0
$ws.Range("Y11").Value = 0
$ws.Range("AI11").Value = 1
$ws.Range("AJ11").Value = 0
$ws.Range("AK11").Value = 0
$ws.Range("AU11").Value = 0.667
$ws.Range("BA11").Value = 2.001
$ws.Range("BM11").Value = 0.667
$ws.Range("BP11").Value = 0.667
$ws.Range("BQ11").Value = 0.7

# Row 12: updated stats
$ws.Range("W12").Value = 1.333
$ws.Range("X12").Value = 0.222
$ws.Range("Y12").Value = 0.471
$ws.Range("AV12").Value = 6
$ws.Range("AW12").Value = 2.449
$ws.Range("BA12").Value = 3.5
$ws.Range("BB12").Value = 0.25
$ws.Range("BC12").Value = 0.5
$ws.Range("BP12").Value = 1.167
$ws.Range("BQ12").Value = 1.208

# Row 13: updated stats
$ws.Range("BP13").Value = 0.722
$ws.Range("BQ13").Value = 0.625
